# Insert a new weekly price row at row 542, pushing the existing rows
# (542-569) down by one. The sheet grows from A1:R569 to A1:R570.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(542).Insert()

$ws.Cells.Item(542, 1).Value = 4
$ws.Cells.Item(542, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(542, 3).Value = "Los Lagos"
$ws.Cells.Item(542, 4).Value = 45267
$ws.Cells.Item(542, 5).Value = 10
$ws.Cells.Item(542, 6).Value = 100112003
$ws.Cells.Item(542, 7).Value = "Ajo"
$ws.Cells.Item(542, 8).Value = "Chino"
$ws.Cells.Item(542, 9).Value = "Primera"
$ws.Cells.Item(542, 10).Value = 120
$ws.Cells.Item(542, 11).Value = 26000
$ws.Cells.Item(542, 12).Value = 26000
$ws.Cells.Item(542, 13).Value = 26000
$ws.Cells.Item(542, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(542, 15).Value = "China"
$ws.Cells.Item(542, 16).Value = 2600
$ws.Cells.Item(542, 17).Value = 10
$ws.Cells.Item(542, 18).Value = "Hortaliza"
